{"js": "// Split two run-on \"Programa resumido\" / \"Crit\u00e9rio\" paragraphs into\n// multiple <w:t> runs separated by manual line breaks (<w:br/>), matching\n// the numbered items / formula pieces that were run together in one string.\n//\n// Word.Range.insertText() accepts the vertical-tab character (\\u000B) as a\n// literal \"manual line break\" marker inline in the replacement text - Word\n// (and this runtime) serialize each \\u000B as a separate <w:t> run boundary\n// joined by <w:br/>, exactly mirroring what a user gets by pressing\n// Shift+Enter between the pieces of text.\n\nconst body = context.document.body;\n\n// --- Change 1: \"Programa resumido\" summary paragraph -----------------\nconst summaryOriginal =\n  \"1- Introdu\u00e7\u00e3o:2- Coeficiente de difus\u00e3o:3- Concentra\u00e7\u00f5es, velocidade e fluxos:\" +\n  \"4 -Equa\u00e7\u00f5es da continuidade em transfer\u00eancia de massa:5- Difus\u00e3o em regime \" +\n  \"permanente sem rea\u00e7\u00e3o qu\u00edmica:6- Difus\u00e3o com rea\u00e7\u00e3o qu\u00edmica:7- Transfer\u00eancia \" +\n  \"de massa entre fases.\";\n\nconst summaryParts = [\n  \"1- Introdu\u00e7\u00e3o:\",\n  \"2- Coeficiente de difus\u00e3o:\",\n  \"3- Concentra\u00e7\u00f5es, velocidade e fluxos:\",\n  \"4 -Equa\u00e7\u00f5es da continuidade em transfer\u00eancia de massa:\",\n  \"5- Difus\u00e3o em regime permanente sem rea\u00e7\u00e3o qu\u00edmica:\",\n  \"6- Difus\u00e3o com rea\u00e7\u00e3o qu\u00edmica:\",\n  \"7- Transfer\u00eancia de massa entre fases.\",\n];\n\nconst summaryResults = body.search(summaryOriginal, { matchCase: true });\nsummaryResults.load(\"items\");\nawait context.sync();\n\nif (summaryResults.items.length > 0) {\n  summaryResults.items[0].insertText(summaryParts.join(\"\\u000B\"), Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Change 2: \"Crit\u00e9rio\" grading-formula sentence --------------------\nconst criterioOriginal =\n  \"A Nota Final (NF) ser\u00e1 calculada da seguinte maneira: NF = (P1 + 2*P2)/3P2 = \" +\n  \"Nota da Prova (80%) e Nota do Trabalho (20%).\";\n\nconst criterioPart1 = \"A Nota Final (NF) ser\u00e1 calculada da seguinte maneira: NF = (P1 + 2*P2)/3\";\nconst criterioPart2 = \"P2 = Nota da Prova (80%) e Nota do Trabalho (20%).\";\n\nconst criterioResults = body.search(criterioOriginal, { matchCase: true });\ncriterioResults.load(\"items\");\nawait context.sync();\n\nif (criterioResults.items.length > 0) {\n  criterioResults.items[0].insertText(\n    criterioPart1 + \"\\u000B\" + criterioPart2,\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n", "ps1": "# Split two run-on \"Programa resumido\" / \"Crit\u00e9rio\" paragraphs into\n# multiple <w:t> runs separated by manual line breaks (<w:br/>), matching\n# the numbered items / formula pieces that were run together in one string.\n#\n# Find.Execute's own \"Replace With\" argument rejects the \"^b\" manual-line-\n# break group character, so instead we locate the target text with\n# Find.Execute (no replace) and then overwrite Range.Text directly with a\n# string containing literal [char]11 (vertical tab / manual line break)\n# markers - Word serializes each one as a separate <w:t> run boundary\n# joined by <w:br/>, exactly mirroring Shift+Enter between the pieces.\n\n$d = $word.ActiveDocument\n$lb = [char]11\n\n# --- Change 1: \"Programa resumido\" summary paragraph -----------------\n$summaryRange = $d.Content\n$summaryNeedle = \"1- Introdu\u00e7\u00e3o:2- Coeficiente de difus\u00e3o:3- Concentra\u00e7\u00f5es, velocidade e fluxos:4 -Equa\u00e7\u00f5es da continuidade em transfer\u00eancia de massa:5- Difus\u00e3o em regime permanente sem rea\u00e7\u00e3o qu\u00edmica:6- Difus\u00e3o com rea\u00e7\u00e3o qu\u00edmica:7- Transfer\u00eancia de massa entre fases.\"\n\n$summaryFound = $summaryRange.Find.Execute($summaryNeedle, $true)\nif ($summaryFound) {\n    $summaryParts = @(\n        \"1- Introdu\u00e7\u00e3o:\",\n        \"2- Coeficiente de difus\u00e3o:\",\n        \"3- Concentra\u00e7\u00f5es, velocidade e fluxos:\",\n        \"4 -Equa\u00e7\u00f5es da continuidade em transfer\u00eancia de massa:\",\n        \"5- Difus\u00e3o em regime permanente sem rea\u00e7\u00e3o qu\u00edmica:\",\n        \"6- Difus\u00e3o com rea\u00e7\u00e3o qu\u00edmica:\",\n        \"7- Transfer\u00eancia de massa entre fases.\"\n    )\n    $summaryRange.Text = [string]::Join($lb, $summaryParts)\n}\n\n# --- Change 2: \"Crit\u00e9rio\" grading-formula sentence --------------------\n$criterioRange = $d.Content\n$criterioNeedle = \"A Nota Final (NF) ser\u00e1 calculada da seguinte maneira: NF = (P1 + 2*P2)/3P2 = Nota da Prova (80%) e Nota do Trabalho (20%).\"\n\n$criterioFound = $criterioRange.Find.Execute($criterioNeedle, $true)\nif ($criterioFound) {\n    $criterioPart1 = \"A Nota Final (NF) ser\u00e1 calculada da seguinte maneira: NF = (P1 + 2*P2)/3\"\n    $criterioPart2 = \"P2 = Nota da Prova (80%) e Nota do Trabalho (20%).\"\n    $criterioRange.Text = $criterioPart1 + $lb + $criterioPart2\n}\n"}
